$wb = $excel.ActiveWorkbook

$wsV003 = $wb.Worksheets.Item("V-003")
$wsCalc = $wb.Worksheets.Item("Calculator")

# --- Calculator sheet: Termin Pembayaran perbulan ---
# Termin Pembayaran (B4) changes from 90 days to 30 days (paid monthly now)
$wsCalc.Activate()
$wsCalc.Range("B4").Value = 30
$wsCalc.Range("E5").Select()

# --- V-003 sheet: Kwitansi Angsuran ada bulan dalam format string ---
# Row 8 (item 5, "Termin Pembayaran perbulan") is now marked "Done" in the Status column (C)
$wsV003.Activate()
$wsV003.Range("C8").Value = "Done"

# Row 12 (item 9, "Calculator Anggsuran") - remove the "Check lagi" note in column D
$wsV003.Range("D12").ClearContents()

# Update selection on V-003 sheet (final active selection in the workbook)
$wsV003.Range("B12").Select()
